$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-OrderRow($r, $name, $phone, $med, $qty, $addr) {
    $ws.Cells.Item($r, 1).Value = $name
    # Phone numbers are stored as text (they look numeric but are not
    # meant to be treated as numbers), so force a text number format
    # before writing the value.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $phone
    $ws.Cells.Item($r, 3).Value = $med
    $ws.Cells.Item($r, 4).Value = $qty
    $ws.Cells.Item($r, 5).Value = $addr
}

# Row 4: replace Pranav's aspirin order with a Chirayu Sahu paracetamol order
Set-OrderRow 4 "Chirayu Sahu" "3234523452" "paracetamol" 1 "Vit Vellore"

# Row 5: update the existing Chirayu Sahu aspirin order to paracetamol
Set-OrderRow 5 "Chirayu Sahu" "3234523452" "paracetamol" 1 "Vit Vellore"

# New rows 6-10: repeated Chirayu Sahu paracetamol orders
foreach ($r in 6..10) {
    Set-OrderRow $r "Chirayu Sahu" "3234523452" "paracetamol" 1 "Vit Vellore"
}

# New rows 11-12: Jon Stewart Doe paracetamol orders
foreach ($r in 11..12) {
    Set-OrderRow $r "Jon Stewart Doe" "6019521325" "paracetamol" 1 "1600 Fake Street"
}
